# Applies the permutation/update of rows 2-10 on the 'viterbi' sheet.
# Each row's tag (column A) changes, and along with it the row's data
# values (columns B:I) move to match the new tag -- i.e. each
# (label, value-row) pair is relocated to a new row index, while the
# underlying data for a given label stays identical.
#
# Numeric literals in scientific notation (e.g. 1e-11) are not parsed
# directly by this interpreter as number literals, so they are built
# from strings and cast to [double].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Num($s) {
    return [double]$s
}

# Final desired state: row number -> (label, values for columns B..I)
$data = @{
    2  = @("NCFS000", 0, 0, 0, 0, 0, 0, 0, 0)
    3  = @("VMIP3S0", (Num "0.0012941074971961"), 0, 0, 0, 0, 0, 0, 0)
    4  = @("AQ0CS0", 0, 0, 0, 0, (Num "6.385259012197449e-11"), 0, 0, 0)
    5  = @("DA0MS0", 0, 0, (Num "2.779888150251198e-06"), 0, 0, 0, 0, 0)
    6  = @("Fp", 0, 0, 0, 0, 0, 0, 0, (Num "3.529518322132861e-16"))
    7  = @("SPS00", 0, (Num "1.988063133517315e-05"), 0, 0, 0, (Num "5.440257583147429e-12"), 0, 0)
    8  = @("AQ0MS0", 0, 0, 0, (Num "1.295582484193506e-09"), 0, 0, 0, 0)
    9  = @("NCMP000", 0, 0, 0, 0, 0, 0, (Num "4.107075865754602e-15"), 0)
    10 = @("NCMS000", 0, 0, 0, (Num "2.827818374535177e-08"), 0, 0, 0, 0)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $row[0]
    for ($c = 2; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
